$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.990.95'
$ws.Range('E2').Value = '  +1.98%  '
$ws.Range('D3').Value = '1.703.92'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.41'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4043'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.474'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '53.26'
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.002'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08824'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '26.09'
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.484'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001356'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.992'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = '1.760.86'
$ws.Range('E17').Value = '  +4.52%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '96.09'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07210'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.335'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.35'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '24.980.22'
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.62'
$ws.Range('E27').Value = '  +3.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.160'
$ws.Range('E28').Value = '  +14.40%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '163.01'
$ws.Range('E29').Value = '  -2.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '151.39'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.356'
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.619'
$ws.Range('E32').Value = '  +17.28%  '
$ws.Range('D33').Value = '1.953.93'
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('E34').Value = '  +6.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08566'
$ws.Range('E35').Value = '  -2.25%  '
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.046'
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2895'
$ws.Range('E38').Value = '  +3.83%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '11.08'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09584'
$ws.Range('E40').Value = '  +4.45%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8325'
$ws.Range('E41').Value = '  +3.22%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '14.06'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.482'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '17.26'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.701'
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.7409'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.257'
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.405'
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.08806'
$ws.Range('E49').Value = '  +7.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.002'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  -0.10%  '

Write-Host "Applied cryptos update"
